# Updated cryptos list on Sat Sep 14 20:59:04 UTC 2024 with GitHub Actions
#
# Applies the latest Price (column D) and Volume(1h) (column E) figures
# to the cryptos worksheet. Values are written as plain text (matching
# the workbook's existing inlineStr/shared-string cells) so that
# numeric-looking prices (e.g. "552.58") are not silently converted into
# real Excel numbers and lose their original textual formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws 'D2' '60.008.09'
Set-TextValue $ws 'E2' '  +0.12%  '

# Row 3 - Ethereum
Set-TextValue $ws 'D3' '2.417.73'
Set-TextValue $ws 'E3' '  -0.26%  '

# Row 4 - TetherUSD
Set-TextValue $ws 'E4' '  -0.01%  '

# Row 5 - BNB
Set-TextValue $ws 'D5' '552.58'
Set-TextValue $ws 'E5' '  -0.37%  '

# Row 6 - Solana
Set-TextValue $ws 'D6' '137.14'
Set-TextValue $ws 'E6' '  -0.93%  '

# Row 7 - USDC
Set-TextValue $ws 'E7' '  +0.01%  '

# Row 8 - XRP
Set-TextValue $ws 'D8' '0.592'

# Row 9 - Dogecoin
Set-TextValue $ws 'E9' '  -1.72%  '

# Row 10 - Toncoin
Set-TextValue $ws 'E10' '  -2.19%  '

# Row 11 - TRON
Set-TextValue $ws 'E11' '  -0.83%  '

# Row 12 - Cardano
Set-TextValue $ws 'D12' '0.354'
Set-TextValue $ws 'E12' '  -1.66%  '

# Row 13 - Avalanche
Set-TextValue $ws 'D13' '25.36'
Set-TextValue $ws 'E13' '  +2.45%  '

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws 'D14' '2.845.88'
Set-TextValue $ws 'E14' '  -0.21%  '

# Row 15 - WrappedBTC
Set-TextValue $ws 'D15' '59.865.40'
Set-TextValue $ws 'E15' '  +0.07%  '

# Row 16 - ShibaInu
Set-TextValue $ws 'E16' '  -1.56%  '

# Row 17 - WrappedEther
Set-TextValue $ws 'D17' '2.428.71'
Set-TextValue $ws 'E17' '  -0.76%  '

# Row 18 - Chainlink
Set-TextValue $ws 'D18' '11.33'
Set-TextValue $ws 'E18' '  -0.82%  '

# Row 19 - Polkadot
Set-TextValue $ws 'E19' '  -0.03%  '

# Row 20 - BitcoinCash
Set-TextValue $ws 'D20' '329.04'
Set-TextValue $ws 'E20' '  -1.58%  '

# Row 21 - Uniswap
Set-TextValue $ws 'D21' '6.68'
Set-TextValue $ws 'E21' '  -3.22%  '

# Row 22 - Dai
Set-TextValue $ws 'E22' '  +0.09%  '

# Row 23 - Litecoin
Set-TextValue $ws 'D23' '65.84'
Set-TextValue $ws 'E23' '  +2.00%  '

# Row 24 - Kaspa
Set-TextValue $ws 'E24' '  +3.54%  '

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws 'E25' '  +0.82%  '

# Row 26 - Binance-PegBSC-USD
Set-TextValue $ws 'E26' '  +0.15%  '

# Row 27 - Fetch.AI
Set-TextValue $ws 'E27' '  +0.00%  '

# Row 28 - PEPE
Set-TextValue $ws 'D28' '0.0₃0776'
Set-TextValue $ws 'E28' '  -1.71%  '

# Row 29 - PancakeSwap
Set-TextValue $ws 'E29' '  -2.10%  '

# Row 30 - Monero
Set-TextValue $ws 'D30' '169.02'
Set-TextValue $ws 'E30' '  -0.90%  '

# Row 31 - Aptos
Set-TextValue $ws 'E31' '  -3.50%  '

# Row 32 - EthereumClassic
Set-TextValue $ws 'D32' '18.62'
Set-TextValue $ws 'E32' '  -0.58%  '

# Row 33 - SuiNetwork
Set-TextValue $ws 'E33' '  -0.25%  '

# Row 34 - USDe
Set-TextValue $ws 'E34' '  +0.00%  '

# Row 36 - FirstDigitalUSD
Set-TextValue $ws 'E36' '  +0.03%  '

# Row 37 - NEARProtocol
Set-TextValue $ws 'E37' '  -1.54%  '

# Row 38 - Stacks
Set-TextValue $ws 'D38' '1.61'
Set-TextValue $ws 'E38' '  -1.84%  '

# Row 39 - Bittensor
Set-TextValue $ws 'D39' '320.55'
Set-TextValue $ws 'E39' '  +2.21%  '

# Row 40 - PolygonEcosystemToken
Set-TextValue $ws 'E40' '  -4.63%  '

# Row 41 - Filecoin
Set-TextValue $ws 'E41' '  -1.82%  '

# Row 42 - Aave
Set-TextValue $ws 'D42' '140.36'
Set-TextValue $ws 'E42' '  -1.58%  '

# Row 43 - Stellar
Set-TextValue $ws 'D43' '0.0969'
Set-TextValue $ws 'E43' '  +0.75%  '

# Row 44 - InjectiveProtocol
Set-TextValue $ws 'D44' '19.56'
Set-TextValue $ws 'E44' '  +1.92%  '

# Row 45 - Hedera
Set-TextValue $ws 'E45' '  -1.76%  '

# Row 46 - Mantle
Set-TextValue $ws 'E46' '  +0.78%  '

# Row 47 - VeChain
Set-TextValue $ws 'E47' '  -1.34%  '

# Row 48 - Polygon
Set-TextValue $ws 'E48' '  -9.01%  '

# Row 49 - WhiteBITCoin
Set-TextValue $ws 'E49' '  +0.01%  '

# Row 50 - dogwifhat
Set-TextValue $ws 'D50' '1.57'
Set-TextValue $ws 'E50' '  -2.78%  '

# Row 51 - ZEEBU
Set-TextValue $ws 'E51' '  -0.90%  '
